$wb = $excel.ActiveWorkbook
Write-Output ($wb | Get-Member | Select-Object -First 60 | Out-String)
